$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 183.22223  # Distill, My Heart
$ws.Range("I9").Value = 196.75  # Distill, My Heart
$ws.Range("K9").Value = 196.75  # Distill, My Heart
$ws.Range("M9").Value = -27.75  # Distill, My Heart
$ws.Range("H19").Value = 987.1053000000001  # Unbreak My Heart
$ws.Range("I19").Value = 1071.7273  # Unbreak My Heart
$ws.Range("J19").Value = 870.75  # Unbreak My Heart
$ws.Range("K19").Value = 1071.7273  # Unbreak My Heart
$ws.Range("L19").Value = 870.75  # Unbreak My Heart
$ws.Range("M19").Value = -896.7273  # Unbreak My Heart
$ws.Range("N19").Value = -1220.75  # Unbreak My Heart
$ws.Range("H137").Value = 4386.227  # Cutting Edge of Culinary Quality
$ws.Range("I137").Value = 3141.5715  # Cutting Edge of Culinary Quality
$ws.Range("J137").Value = 4967.067  # Cutting Edge of Culinary Quality
$ws.Range("K137").Value = 9424.7145  # Cutting Edge of Culinary Quality
$ws.Range("L137").Value = 14901.201  # Cutting Edge of Culinary Quality
$ws.Range("M137").Value = -6874.7145  # Cutting Edge of Culinary Quality
$ws.Range("N137").Value = -20001.201  # Cutting Edge of Culinary Quality
$ws.Range("H138").Value = 2160.0833  # All-night Crafting
$ws.Range("I138").Value = 1518.5714  # All-night Crafting
$ws.Range("J138").Value = 3058.2  # All-night Crafting
$ws.Range("K138").Value = 4555.7142  # All-night Crafting
$ws.Range("L138").Value = 9174.599999999999  # All-night Crafting
$ws.Range("M138").Value = 584.2857999999997  # All-night Crafting
$ws.Range("N138").Value = -19454.6  # All-night Crafting

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4702.8  # Hollow Hallmarks
$ws.Range("I45").Value = 2750  # Hollow Hallmarks
$ws.Range("J45").Value = 5191  # Hollow Hallmarks
$ws.Range("K45").Value = 2750  # Hollow Hallmarks
$ws.Range("L45").Value = 5191  # Hollow Hallmarks
$ws.Range("M45").Value = -2373  # Hollow Hallmarks
$ws.Range("N45").Value = -5945  # Hollow Hallmarks
$ws.Range("H61").Value = 3799.8235  # Dealing with the Tough Stuff
$ws.Range("I61").Value = 2507.4614  # Dealing with the Tough Stuff
$ws.Range("J61").Value = 8000  # Dealing with the Tough Stuff
$ws.Range("K61").Value = 2507.4614  # Dealing with the Tough Stuff
$ws.Range("L61").Value = 8000  # Dealing with the Tough Stuff
$ws.Range("M61").Value = -2295.4614  # Dealing with the Tough Stuff
$ws.Range("N61").Value = -8424  # Dealing with the Tough Stuff
$ws.Range("H74").Value = 679.8  # As the Bolt Flies
$ws.Range("I74").Value = 599.75  # As the Bolt Flies
$ws.Range("J74").Value = 1000  # As the Bolt Flies
$ws.Range("K74").Value = 599.75  # As the Bolt Flies
$ws.Range("L74").Value = 1000  # As the Bolt Flies
$ws.Range("M74").Value = 274.25  # As the Bolt Flies
$ws.Range("N74").Value = -2748  # As the Bolt Flies
$ws.Range("H77").Value = 679.8  # Heavy Metal Banned (L)
$ws.Range("I77").Value = 599.75  # Heavy Metal Banned (L)
$ws.Range("J77").Value = 1000  # Heavy Metal Banned (L)
$ws.Range("K77").Value = 2998.75  # Heavy Metal Banned (L)
$ws.Range("L77").Value = 5000  # Heavy Metal Banned (L)
$ws.Range("M77").Value = 1369.25  # Heavy Metal Banned (L)
$ws.Range("N77").Value = -13736  # Heavy Metal Banned (L)
$ws.Range("H97").Value = 3584990  # Ore for Me
$ws.Range("J97").Value = 15873672  # Ore for Me
$ws.Range("L97").Value = 15873672  # Ore for Me
$ws.Range("N97").Value = -15874664  # Ore for Me
$ws.Range("H132").Value = 2343  # Don't Bore Me, Ore Me
$ws.Range("I132").Value = 1371.6  # Don't Bore Me, Ore Me
$ws.Range("J132").Value = 7200  # Don't Bore Me, Ore Me
$ws.Range("K132").Value = 4114.799999999999  # Don't Bore Me, Ore Me
$ws.Range("L132").Value = 21600  # Don't Bore Me, Ore Me
$ws.Range("M132").Value = -1584.799999999999  # Don't Bore Me, Ore Me
$ws.Range("N132").Value = -26660  # Don't Bore Me, Ore Me
$ws.Range("H136").Value = 3799.8235  # Metal with Mettle
$ws.Range("I136").Value = 2507.4614  # Metal with Mettle
$ws.Range("J136").Value = 8000  # Metal with Mettle
$ws.Range("K136").Value = 7522.3842  # Metal with Mettle
$ws.Range("L136").Value = 24000  # Metal with Mettle
$ws.Range("M136").Value = -4972.3842  # Metal with Mettle
$ws.Range("N136").Value = -29100  # Metal with Mettle

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 5052429.5  # High Steal
$ws.Range("I94").Value = 1397.591  # High Steal
$ws.Range("K94").Value = 1397.591  # High Steal
$ws.Range("M94").Value = -946.5909999999999  # High Steal
$ws.Range("H134").Value = 2356.2456  # Ruthenium Supremium
$ws.Range("I134").Value = 1547.9791  # Ruthenium Supremium
$ws.Range("J134").Value = 6667  # Ruthenium Supremium
$ws.Range("K134").Value = 4643.9373  # Ruthenium Supremium
$ws.Range("L134").Value = 20001  # Ruthenium Supremium
$ws.Range("M134").Value = -2108.9373  # Ruthenium Supremium
$ws.Range("N134").Value = -25071  # Ruthenium Supremium

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3207.2144  # Wall Not Found
$ws.Range("I31").Value = 1587.5294  # Wall Not Found
$ws.Range("J31").Value = 5710.364  # Wall Not Found
$ws.Range("K31").Value = 1587.5294  # Wall Not Found
$ws.Range("L31").Value = 5710.364  # Wall Not Found
$ws.Range("M31").Value = -1292.5294  # Wall Not Found
$ws.Range("N31").Value = -6300.364  # Wall Not Found
$ws.Range("H34").Value = 3207.2144  # Armoires of the Rich and Famous
$ws.Range("I34").Value = 1587.5294  # Armoires of the Rich and Famous
$ws.Range("J34").Value = 5710.364  # Armoires of the Rich and Famous
$ws.Range("K34").Value = 1587.5294  # Armoires of the Rich and Famous
$ws.Range("L34").Value = 5710.364  # Armoires of the Rich and Famous
$ws.Range("M34").Value = -1385.5294  # Armoires of the Rich and Famous
$ws.Range("N34").Value = -6114.364  # Armoires of the Rich and Famous
$ws.Range("H94").Value = 7277.6875  # Beech, Please
$ws.Range("I94").Value = 20386.6  # Beech, Please
$ws.Range("K94").Value = 20386.6  # Beech, Please
$ws.Range("M94").Value = -19935.6  # Beech, Please
$ws.Range("I99").Value = 2442891  # O Pine
$ws.Range("K99").Value = 2442891  # O Pine
$ws.Range("M99").Value = -2441393  # O Pine
$ws.Range("I126").Value = 2442891  # A Better Conductor
$ws.Range("K126").Value = 7328673  # A Better Conductor
$ws.Range("M126").Value = -7326203  # A Better Conductor

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 98733.69500000001  # Sky Is the Limit
$ws.Range("I70").Value = 113049.55  # Sky Is the Limit
$ws.Range("K70").Value = 113049.55  # Sky Is the Limit
$ws.Range("M70").Value = -112779.55  # Sky Is the Limit
$ws.Range("H73").Value = 98733.69500000001  # Hulls of Broken Dreams (L)
$ws.Range("I73").Value = 113049.55  # Hulls of Broken Dreams (L)
$ws.Range("K73").Value = 113049.55  # Hulls of Broken Dreams (L)
$ws.Range("M73").Value = -112113.55  # Hulls of Broken Dreams (L)
$ws.Range("H102").Value = 2027.3334  # Put the Metal to the Peddle
$ws.Range("I102").Value = 1230.0555  # Put the Metal to the Peddle
$ws.Range("J102").Value = 3621.889  # Put the Metal to the Peddle
$ws.Range("K102").Value = 1230.0555  # Put the Metal to the Peddle
$ws.Range("L102").Value = 3621.889  # Put the Metal to the Peddle
$ws.Range("M102").Value = 391.9445000000001  # Put the Metal to the Peddle
$ws.Range("N102").Value = -6865.889  # Put the Metal to the Peddle
$ws.Range("H113").Value = 2267.75  # Copious Crystal Cannons
$ws.Range("I113").Value = 1197.5714  # Copious Crystal Cannons
$ws.Range("J113").Value = 5478.2856  # Copious Crystal Cannons
$ws.Range("K113").Value = 1197.5714  # Copious Crystal Cannons
$ws.Range("L113").Value = 5478.2856  # Copious Crystal Cannons
$ws.Range("M113").Value = 972.4286  # Copious Crystal Cannons
$ws.Range("N113").Value = -9818.285599999999  # Copious Crystal Cannons
$ws.Range("H126").Value = 5121.8887  # Gold Rush Order
$ws.Range("I126").Value = 3874.5  # Gold Rush Order
$ws.Range("J126").Value = 6119.8  # Gold Rush Order
$ws.Range("K126").Value = 11623.5  # Gold Rush Order
$ws.Range("L126").Value = 18359.4  # Gold Rush Order
$ws.Range("M126").Value = -9153.5  # Gold Rush Order
$ws.Range("N126").Value = -23299.4  # Gold Rush Order
$ws.Range("H132").Value = 2228.9443  # On Board for Lar
$ws.Range("I132").Value = 2065.9412  # On Board for Lar
$ws.Range("K132").Value = 6197.823600000001  # On Board for Lar
$ws.Range("M132").Value = -3667.823600000001  # On Board for Lar

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4924.9165  # Tan Before the Ban
$ws.Range("J7").Value = 5635.1113  # Tan Before the Ban
$ws.Range("L7").Value = 5635.1113  # Tan Before the Ban
$ws.Range("N7").Value = -5859.1113  # Tan Before the Ban
$ws.Range("H122").Value = 5076.0713  # Hell on Leather
$ws.Range("I122").Value = 2558.5715  # Hell on Leather
$ws.Range("J122").Value = 7593.5713  # Hell on Leather
$ws.Range("K122").Value = 7675.7145  # Hell on Leather
$ws.Range("L122").Value = 22780.7139  # Hell on Leather
$ws.Range("M122").Value = -5225.7145  # Hell on Leather
$ws.Range("N122").Value = -27680.7139  # Hell on Leather
$ws.Range("H126").Value = 4924.9165  # Battered Books
$ws.Range("J126").Value = 5635.1113  # Battered Books
$ws.Range("L126").Value = 16905.3339  # Battered Books
$ws.Range("N126").Value = -21845.3339  # Battered Books
$ws.Range("H136").Value = 3651.5264  # Respect for Br'aax
$ws.Range("I136").Value = 1940.6923  # Respect for Br'aax
$ws.Range("J136").Value = 7358.3335  # Respect for Br'aax
$ws.Range("K136").Value = 5822.0769  # Respect for Br'aax
$ws.Range("L136").Value = 22075.0005  # Respect for Br'aax
$ws.Range("M136").Value = -3272.0769  # Respect for Br'aax
$ws.Range("N136").Value = -27175.0005  # Respect for Br'aax

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 23749  # After the Smock-down
$ws.Range("I51").Value = 20000  # After the Smock-down
$ws.Range("J51").Value = 27498  # After the Smock-down
$ws.Range("K51").Value = 20000  # After the Smock-down
$ws.Range("L51").Value = 27498  # After the Smock-down
$ws.Range("M51").Value = -19490  # After the Smock-down
$ws.Range("N51").Value = -28518  # After the Smock-down
$ws.Range("H81").Value = 15153839  # Where the Dragonflies, the Net Catches
$ws.Range("J81").Value = 30304968  # Where the Dragonflies, the Net Catches
$ws.Range("L81").Value = 60609936  # Where the Dragonflies, the Net Catches
$ws.Range("N81").Value = -60612058  # Where the Dragonflies, the Net Catches
$ws.Range("H84").Value = 15153839  # To Kill a Dragon on Nameday (L)
$ws.Range("J84").Value = 30304968  # To Kill a Dragon on Nameday (L)
$ws.Range("L84").Value = 303049680  # To Kill a Dragon on Nameday (L)
$ws.Range("N84").Value = -303060288  # To Kill a Dragon on Nameday (L)
$ws.Range("H132").Value = 3084.238  # Comfy Cabins
$ws.Range("I132").Value = 1827  # Comfy Cabins
$ws.Range("K132").Value = 5481  # Comfy Cabins
$ws.Range("M132").Value = -2951  # Comfy Cabins
$ws.Range("H136").Value = 2033.5625  # Weaving the Envelope
$ws.Range("I136").Value = 1095.4445  # Weaving the Envelope
$ws.Range("J136").Value = 3239.7144  # Weaving the Envelope
$ws.Range("K136").Value = 3286.3335  # Weaving the Envelope
$ws.Range("L136").Value = 9719.143199999999  # Weaving the Envelope
$ws.Range("M136").Value = -736.3335000000002  # Weaving the Envelope
$ws.Range("N136").Value = -14819.1432  # Weaving the Envelope
